# Apply the "Updated cryptos list" refresh: new Price (D) / Volume(1h) (E)
# values per row, plus rows 48-49 swapping ARBITRUM <-> Filecoin (with their
# own refreshed price/volume), matching the upstream GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value2 = "72.449.47"
$ws.Cells.Item(2,5).Value2 = "  +0.19%  "
$ws.Cells.Item(3,4).Value2 = "2.647.14"
$ws.Cells.Item(3,5).Value2 = "  -0.92%  "
$ws.Cells.Item(4,5).Value2 = "  -0.03%  "
$ws.Cells.Item(5,4).Value2 = "588.93"
$ws.Cells.Item(5,5).Value2 = "  -2.26%  "
$ws.Cells.Item(6,4).Value2 = "174.75"
$ws.Cells.Item(6,5).Value2 = "  -2.28%  "
$ws.Cells.Item(7,5).Value2 = "  -0.02%  "
$ws.Cells.Item(8,5).Value2 = "  -1.04%  "
$ws.Cells.Item(9,4).Value2 = "0.172"
$ws.Cells.Item(9,5).Value2 = "  -0.79%  "
$ws.Cells.Item(10,4).Value2 = "2.646.70"
$ws.Cells.Item(10,5).Value2 = "  -0.93%  "
$ws.Cells.Item(11,5).Value2 = "  +1.04%  "
$ws.Cells.Item(12,4).Value2 = "0.356"
$ws.Cells.Item(12,5).Value2 = "  -0.41%  "
$ws.Cells.Item(13,4).Value2 = "4.95"
$ws.Cells.Item(13,5).Value2 = "  -1.42%  "
$ws.Cells.Item(14,4).Value2 = "3.130.09"
$ws.Cells.Item(14,5).Value2 = "  -0.88%  "
$ws.Cells.Item(15,4).Value2 = "0.0000187"
$ws.Cells.Item(15,5).Value2 = "  -1.26%  "
$ws.Cells.Item(16,4).Value2 = "72.255.29"
$ws.Cells.Item(16,5).Value2 = "  +0.05%  "
$ws.Cells.Item(17,4).Value2 = "25.95"
$ws.Cells.Item(17,5).Value2 = "  -2.59%  "
$ws.Cells.Item(18,4).Value2 = "2.637.35"
$ws.Cells.Item(18,5).Value2 = "  -0.20%  "
$ws.Cells.Item(19,4).Value2 = "12.15"
$ws.Cells.Item(19,5).Value2 = "  +1.21%  "
$ws.Cells.Item(20,4).Value2 = "7.99"
$ws.Cells.Item(20,5).Value2 = "  -0.41%  "
$ws.Cells.Item(21,4).Value2 = "371.83"
$ws.Cells.Item(21,5).Value2 = "  -2.01%  "
$ws.Cells.Item(22,4).Value2 = "4.16"
$ws.Cells.Item(22,5).Value2 = "  -1.04%  "
$ws.Cells.Item(23,4).Value2 = "2.06"
$ws.Cells.Item(23,5).Value2 = "  -0.29%  "
$ws.Cells.Item(24,5).Value2 = "  -0.01%  "
$ws.Cells.Item(25,4).Value2 = "71.00"
$ws.Cells.Item(25,5).Value2 = "  -2.11%  "
$ws.Cells.Item(26,4).Value2 = "4.26"
$ws.Cells.Item(26,5).Value2 = "  -3.16%  "
$ws.Cells.Item(27,4).Value2 = "9.66"
$ws.Cells.Item(27,5).Value2 = "  -3.59%  "
$ws.Cells.Item(28,4).Value2 = "2.781.54"
$ws.Cells.Item(28,5).Value2 = "  -0.93%  "
$ws.Cells.Item(29,5).Value2 = "  -0.36%  "
$ws.Cells.Item(30,4).Value2 = "0.0₃0957"
$ws.Cells.Item(30,5).Value2 = "  +0.67%  "
$ws.Cells.Item(31,4).Value2 = "8.02"
$ws.Cells.Item(31,5).Value2 = "  -2.35%  "
$ws.Cells.Item(32,4).Value2 = "498.05"
$ws.Cells.Item(32,5).Value2 = "  -4.83%  "
$ws.Cells.Item(33,5).Value2 = "  -2.32%  "
$ws.Cells.Item(34,5).Value2 = "  -1.26%  "
$ws.Cells.Item(35,5).Value2 = "  -0.01%  "
$ws.Cells.Item(36,4).Value2 = "162.50"
$ws.Cells.Item(36,5).Value2 = "  -1.41%  "
$ws.Cells.Item(37,4).Value2 = "19.28"
$ws.Cells.Item(37,5).Value2 = "  -1.66%  "
$ws.Cells.Item(38,4).Value2 = "0.114"
$ws.Cells.Item(38,5).Value2 = "  +2.15%  "
$ws.Cells.Item(39,4).Value2 = "18.88"
$ws.Cells.Item(39,5).Value2 = "  -1.26%  "
$ws.Cells.Item(40,5).Value2 = "  -2.53%  "
$ws.Cells.Item(41,5).Value2 = "  -0.06%  "
$ws.Cells.Item(42,4).Value2 = "1.74"
$ws.Cells.Item(42,5).Value2 = "  -5.78%  "
$ws.Cells.Item(43,4).Value2 = "2.58"
$ws.Cells.Item(43,5).Value2 = "  -0.87%  "
$ws.Cells.Item(44,4).Value2 = "4.91"
$ws.Cells.Item(44,5).Value2 = "  -3.22%  "
$ws.Cells.Item(45,4).Value2 = "0.328"
$ws.Cells.Item(45,5).Value2 = "  -2.35%  "
$ws.Cells.Item(46,5).Value2 = "  -0.68%  "
$ws.Cells.Item(47,4).Value2 = "152.90"
$ws.Cells.Item(47,5).Value2 = "  -0.29%  "
$ws.Cells.Item(50,4).Value2 = "1.68"
$ws.Cells.Item(50,5).Value2 = "  -2.69%  "
$ws.Cells.Item(51,4).Value2 = "0.0748"
$ws.Cells.Item(51,5).Value2 = "  -2.09%  "

# Row 48: was Filecoin -> now ARBITRUM (moved up, refreshed price/volume)
$ws.Cells.Item(48,2).Value2 = "ARBITRUM"
$ws.Cells.Item(48,3).Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(48,4).Value2 = "0.549"
$ws.Cells.Item(48,5).Value2 = "  -0.48%  "

# Row 49: was ARBITRUM -> now Filecoin (moved down, refreshed price/volume)
$ws.Cells.Item(49,2).Value2 = "Filecoin"
$ws.Cells.Item(49,3).Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(49,4).Value2 = "3.66"
$ws.Cells.Item(49,5).Value2 = "  -2.62%  "
